$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("multi_shop_failures")

$ws.Range("B2").Value = "2022 Jeep Grand Cherokee Enhanced Powertrain"
$ws.Range("B3").Value = "2022 Jeep Grand Cherokee Enhanced Powertrain"
$ws.Range("B4").Value = "2021 Jeep Grand Cherokee Enhanced Powertrain"
$ws.Range("B5").Value = "2021 Jeep Grand Cherokee Enhanced Powertrain"
$ws.Range("B6").Value = "2021 Jeep Grand Cherokee Enhanced Powertrain"
$ws.Range("B7").Value = "2024 Jeep Grand Wagoneer L Enhanced Powertrain"
$ws.Range("B8").Value = "2024 Jeep Wagoneer Enhanced Powertrain"
$ws.Range("B9").Value = "2024 Jeep Wagoneer Enhanced Powertrain"
$ws.Range("B10").Value = "2024 Dodge Ram 1500 Enhanced Powertrain"
$ws.Range("B11").Value = "2023 Dodge Ram 1500 Enhanced Powertrain"
$ws.Range("B12").Value = "2017 Jeep Compass Enhanced Powertrain"
$ws.Range("B13").Value = "2024 Chevrolet Blazer Enhanced Powertrain"
$ws.Range("B14").Value = "2017 Kia Forte 4 Door Enhanced Powertrain UDS - Theta GDI 1.6 / 2.0L"
$ws.Range("B15").Value = "2020 Toyota RAV4 Enhanced Powertrain CAN"
$ws.Range("B16").Value = "2006 Subaru Outback Enhanced Powertrain"
$ws.Range("B17").Value = "2018 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B18").Value = "2017 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B19").Value = "2018 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B20").Value = "2018 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B21").Value = "2017 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B22").Value = "2017 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B23").Value = "2017 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B24").Value = "2017 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B25").Value = "2017 Hyundai Santa Fe 2.4L Enhanced Powertrain UDS - Theta 2.XL"
$ws.Range("B26").Value = "2020 Volvo S60 Engine Control Module (ECM)"
$ws.Range("B27").Value = "2024 Mitsubishi Outlander PHEV Enhanced Powertrain CAN"
$ws.Range("B28").Value = "2024 Mitsubishi Outlander PHEV Enhanced Powertrain CAN"
$ws.Range("B29").Value = "2024 Mitsubishi Outlander PHEV Enhanced Powertrain CAN"
$ws.Range("B30").Value = "2022 Subaru Impreza Enhanced Powertrain CAN"
$ws.Range("B31").Value = "2023 Subaru Impreza Enhanced Powertrain CAN"
$ws.Range("B32").Value = "2024 Alfa-Romeo Stelvio Enhanced Powertrain"
$ws.Range("B33").Value = "2019 Toyota Corolla Enhanced Powertrain CAN"
$ws.Range("B34").Value = "2022 Hyundai Kona Enhanced Powertrain UDS - T-GDI 1.6L Gamma"
$ws.Range("B35").Value = "2023 Hyundai Kona Enhanced Powertrain UDS - T-GDI 1.6L Gamma"
$ws.Range("B36").Value = "2022 Hyundai Kona Enhanced Powertrain UDS - MPI 2.0L"
$ws.Range("B37").Value = "2023 Hyundai Kona Enhanced Powertrain UDS - T-GDI 1.6L Gamma"
$ws.Range("B38").Value = "2019 Volvo XC60 Engine Control Module (ECM)"
$ws.Range("B39").Value = "2024 Mercedes-Benz C Class Enhanced Powertrain CAN MED41"
$ws.Range("B40").Value = "2007 BMW 3 Series Enhanced Powertrain"
$ws.Range("B41").Value = "2009 Volvo S40 Enhanced Powertrain CAN"
$ws.Range("B42").Value = "2023 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B43").Value = "2019 Volvo XC40 Supplemental Restraint System Module (SRS)"
$ws.Range("B44").Value = "2018 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B45").Value = "2019 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B46").Value = "2017 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B47").Value = "2019 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B48").Value = "2016 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B49").Value = "2022 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B50").Value = "2017 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B51").Value = "2018 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B52").Value = "2019 Volvo XC90 Engine Control Module (ECM)"
$ws.Range("B53").Value = "2018 Volvo XC60 Supplemental Restraint System Module (SRS)"
$ws.Range("B54").Value = "2021 Volvo XC60 Engine Control Module (ECM)"
$ws.Range("B55").Value = "2019 Dodge Ram 1500 Enhanced Powertrain"
